$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (col C) and montant_total (col D) values for the
# 2020-08-20 Fonds de solidarite data refresh.

$ws.Cells.Item(2, 3).Value = 38120
$ws.Cells.Item(2, 4).Value = 55130766
$ws.Cells.Item(3, 3).Value = 91686
$ws.Cells.Item(3, 4).Value = 134397491
$ws.Cells.Item(4, 3).Value = 31371
$ws.Cells.Item(4, 4).Value = 46461119
$ws.Cells.Item(5, 3).Value = 8762
$ws.Cells.Item(5, 4).Value = 13023563
$ws.Cells.Item(6, 3).Value = 2019
$ws.Cells.Item(6, 4).Value = 3000471
$ws.Cells.Item(7, 3).Value = 156
$ws.Cells.Item(7, 4).Value = 229093
$ws.Cells.Item(12, 3).Value = 41616
$ws.Cells.Item(12, 4).Value = 56469057
$ws.Cells.Item(13, 3).Value = 9752
$ws.Cells.Item(13, 4).Value = 14107440
$ws.Cells.Item(14, 3).Value = 26132
$ws.Cells.Item(14, 4).Value = 38321573
$ws.Cells.Item(15, 3).Value = 8360
$ws.Cells.Item(15, 4).Value = 12407824
$ws.Cells.Item(16, 3).Value = 2163
$ws.Cells.Item(16, 4).Value = 3216665
$ws.Cells.Item(17, 3).Value = 423
$ws.Cells.Item(17, 4).Value = 623623
$ws.Cells.Item(20, 3).Value = 10296
$ws.Cells.Item(20, 4).Value = 13624667
$ws.Cells.Item(21, 3).Value = 13519
$ws.Cells.Item(21, 4).Value = 19515095
$ws.Cells.Item(22, 3).Value = 31846
$ws.Cells.Item(22, 4).Value = 46731228
$ws.Cells.Item(23, 3).Value = 10278
$ws.Cells.Item(23, 4).Value = 15278578
$ws.Cells.Item(24, 3).Value = 2659
$ws.Cells.Item(24, 4).Value = 3953682
$ws.Cells.Item(25, 3).Value = 511
$ws.Cells.Item(25, 4).Value = 760592
$ws.Cells.Item(27, 3).Value = 11760
$ws.Cells.Item(27, 4).Value = 15705041
$ws.Cells.Item(28, 3).Value = 7741
$ws.Cells.Item(28, 4).Value = 11206295
$ws.Cells.Item(29, 3).Value = 22665
$ws.Cells.Item(29, 4).Value = 33270373
$ws.Cells.Item(30, 3).Value = 7868
$ws.Cells.Item(30, 4).Value = 11706943
$ws.Cells.Item(31, 3).Value = 1984
$ws.Cells.Item(31, 4).Value = 2960499
$ws.Cells.Item(34, 3).Value = 8369
$ws.Cells.Item(34, 4).Value = 11054100
$ws.Cells.Item(35, 3).Value = 3281
$ws.Cells.Item(35, 4).Value = 4737437
$ws.Cells.Item(36, 3).Value = 7916
$ws.Cells.Item(36, 4).Value = 11561227
$ws.Cells.Item(37, 3).Value = 3197
$ws.Cells.Item(37, 4).Value = 4738461
$ws.Cells.Item(38, 3).Value = 831
$ws.Cells.Item(38, 4).Value = 1237723
$ws.Cells.Item(39, 3).Value = 165
$ws.Cells.Item(39, 4).Value = 245186
$ws.Cells.Item(41, 3).Value = 2497
$ws.Cells.Item(41, 4).Value = 3377423
$ws.Cells.Item(42, 3).Value = 17413
$ws.Cells.Item(42, 4).Value = 25178971
$ws.Cells.Item(43, 3).Value = 51531
$ws.Cells.Item(43, 4).Value = 75542243
$ws.Cells.Item(44, 3).Value = 19115
$ws.Cells.Item(44, 4).Value = 28393508
$ws.Cells.Item(45, 3).Value = 5645
$ws.Cells.Item(45, 4).Value = 8403805
$ws.Cells.Item(46, 3).Value = 1220
$ws.Cells.Item(46, 4).Value = 1820545
$ws.Cells.Item(50, 3).Value = 16834
$ws.Cells.Item(50, 4).Value = 22403075
$ws.Cells.Item(51, 3).Value = 2062
$ws.Cells.Item(51, 4).Value = 2990096
$ws.Cells.Item(52, 3).Value = 7009
$ws.Cells.Item(52, 4).Value = 10302108
$ws.Cells.Item(53, 3).Value = 2379
$ws.Cells.Item(53, 4).Value = 3552964
$ws.Cells.Item(54, 3).Value = 759
$ws.Cells.Item(54, 4).Value = 1133805
$ws.Cells.Item(55, 3).Value = 188
$ws.Cells.Item(55, 4).Value = 278833
$ws.Cells.Item(57, 3).Value = 7114
$ws.Cells.Item(57, 4).Value = 9780005
$ws.Cells.Item(58, 3).Value = 1034
$ws.Cells.Item(58, 4).Value = 1629621
$ws.Cells.Item(59, 3).Value = 2588
$ws.Cells.Item(59, 4).Value = 4105551
$ws.Cells.Item(60, 3).Value = 1029
$ws.Cells.Item(60, 4).Value = 1645539
$ws.Cells.Item(61, 3).Value = 348
$ws.Cells.Item(61, 4).Value = 555883
$ws.Cells.Item(62, 3).Value = 114
$ws.Cells.Item(62, 4).Value = 185850
$ws.Cells.Item(63, 3).Value = 21
$ws.Cells.Item(63, 4).Value = 36000
$ws.Cells.Item(64, 3).Value = 1496
$ws.Cells.Item(64, 4).Value = 2217068
$ws.Cells.Item(65, 3).Value = 15536
$ws.Cells.Item(65, 4).Value = 22436832
$ws.Cells.Item(66, 3).Value = 45068
$ws.Cells.Item(66, 4).Value = 65944967
$ws.Cells.Item(67, 3).Value = 15801
$ws.Cells.Item(67, 4).Value = 23481042
$ws.Cells.Item(68, 3).Value = 4596
$ws.Cells.Item(68, 4).Value = 6845051
$ws.Cells.Item(69, 3).Value = 933
$ws.Cells.Item(69, 4).Value = 1387668
$ws.Cells.Item(73, 3).Value = 15195
$ws.Cells.Item(73, 4).Value = 20023257
$ws.Cells.Item(74, 3).Value = 52604
$ws.Cells.Item(74, 4).Value = 76557273
$ws.Cells.Item(75, 3).Value = 148494
$ws.Cells.Item(75, 4).Value = 218778738
$ws.Cells.Item(76, 3).Value = 64285
$ws.Cells.Item(76, 4).Value = 95795302
$ws.Cells.Item(77, 3).Value = 20548
$ws.Cells.Item(77, 4).Value = 30702322
$ws.Cells.Item(78, 3).Value = 4887
$ws.Cells.Item(78, 4).Value = 7299401
$ws.Cells.Item(85, 3).Value = 51854
$ws.Cells.Item(85, 4).Value = 70540439
$ws.Cells.Item(86, 3).Value = 4668
$ws.Cells.Item(86, 4).Value = 6765489
$ws.Cells.Item(87, 3).Value = 11687
$ws.Cells.Item(87, 4).Value = 17169363
$ws.Cells.Item(88, 3).Value = 3910
$ws.Cells.Item(88, 4).Value = 5827583
$ws.Cells.Item(89, 3).Value = 1353
$ws.Cells.Item(89, 4).Value = 2021989
$ws.Cells.Item(93, 3).Value = 5461
$ws.Cells.Item(93, 4).Value = 7341918
$ws.Cells.Item(94, 3).Value = 1614
$ws.Cells.Item(94, 4).Value = 2324533
$ws.Cells.Item(95, 3).Value = 5234
$ws.Cells.Item(95, 4).Value = 7709692
$ws.Cells.Item(96, 3).Value = 1949
$ws.Cells.Item(96, 4).Value = 2903437
$ws.Cells.Item(97, 3).Value = 698
$ws.Cells.Item(97, 4).Value = 1045960
$ws.Cells.Item(101, 3).Value = 3603
$ws.Cells.Item(101, 4).Value = 4769610
$ws.Cells.Item(102, 3).Value = 656
$ws.Cells.Item(102, 4).Value = 1041709
$ws.Cells.Item(103, 3).Value = 391
$ws.Cells.Item(103, 4).Value = 631592
$ws.Cells.Item(104, 3).Value = 138
$ws.Cells.Item(104, 4).Value = 214160
$ws.Cells.Item(106, 3).Value = 23
$ws.Cells.Item(106, 4).Value = 37500
$ws.Cells.Item(107, 3).Value = 10890
$ws.Cells.Item(107, 4).Value = 15798814
$ws.Cells.Item(108, 3).Value = 29418
$ws.Cells.Item(108, 4).Value = 43217518
$ws.Cells.Item(109, 3).Value = 9847
$ws.Cells.Item(109, 4).Value = 14643179
$ws.Cells.Item(110, 3).Value = 2711
$ws.Cells.Item(110, 4).Value = 4042707
$ws.Cells.Item(114, 3).Value = 9859
$ws.Cells.Item(114, 4).Value = 13023974
$ws.Cells.Item(115, 3).Value = 30761
$ws.Cells.Item(115, 4).Value = 44362400
$ws.Cells.Item(116, 3).Value = 66615
$ws.Cells.Item(116, 4).Value = 97486102
$ws.Cells.Item(117, 3).Value = 21499
$ws.Cells.Item(117, 4).Value = 31950713
$ws.Cells.Item(118, 3).Value = 6101
$ws.Cells.Item(118, 4).Value = 9089521
$ws.Cells.Item(119, 3).Value = 1137
$ws.Cells.Item(119, 4).Value = 1699271
$ws.Cells.Item(120, 3).Value = 80
$ws.Cells.Item(120, 4).Value = 117420
$ws.Cells.Item(124, 3).Value = 26040
$ws.Cells.Item(124, 4).Value = 34774512
$ws.Cells.Item(125, 3).Value = 36375
$ws.Cells.Item(125, 4).Value = 52496517
$ws.Cells.Item(126, 3).Value = 77358
$ws.Cells.Item(126, 4).Value = 113121597
$ws.Cells.Item(127, 3).Value = 24004
$ws.Cells.Item(127, 4).Value = 35626437
$ws.Cells.Item(128, 3).Value = 6431
$ws.Cells.Item(128, 4).Value = 9557238
$ws.Cells.Item(129, 3).Value = 1247
$ws.Cells.Item(129, 4).Value = 1854911
$ws.Cells.Item(133, 3).Value = 32010
$ws.Cells.Item(133, 4).Value = 42498669
$ws.Cells.Item(134, 3).Value = 13384
$ws.Cells.Item(134, 4).Value = 19376293
$ws.Cells.Item(135, 3).Value = 32549
$ws.Cells.Item(135, 4).Value = 47804122
$ws.Cells.Item(136, 3).Value = 11543
$ws.Cells.Item(136, 4).Value = 17151392
$ws.Cells.Item(137, 3).Value = 2975
$ws.Cells.Item(137, 4).Value = 4433741
$ws.Cells.Item(138, 3).Value = 504
$ws.Cells.Item(138, 4).Value = 749990
$ws.Cells.Item(141, 3).Value = 10892
$ws.Cells.Item(141, 4).Value = 14522513
$ws.Cells.Item(142, 3).Value = 35447
$ws.Cells.Item(142, 4).Value = 51195951
$ws.Cells.Item(143, 3).Value = 81970
$ws.Cells.Item(143, 4).Value = 120091654
$ws.Cells.Item(144, 3).Value = 24529
$ws.Cells.Item(144, 4).Value = 36442453
$ws.Cells.Item(145, 3).Value = 6441
$ws.Cells.Item(145, 4).Value = 9611067
$ws.Cells.Item(146, 3).Value = 1452
$ws.Cells.Item(146, 4).Value = 2160230
$ws.Cells.Item(149, 3).Value = 29399
$ws.Cells.Item(149, 4).Value = 39653653
